# Add two new rows at the bottom of the data, continuing the date sequence
# in column A by one day each, and repeating the same B:J values as the
# last existing row (matches the pattern already present in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

# Capture the values currently on the last row (B:lastCol) so the new rows
# replicate them, and capture the date serial in column A to increment from.
$lastDate = $ws.Cells.Item($lastRow, 1).Value2()
$rowValues = @()
for ($c = 2; $c -le $lastCol; $c++) {
    $rowValues += , ($ws.Cells.Item($lastRow, $c).Value2())
}

$rowsToAdd = 2
for ($i = 1; $i -le $rowsToAdd; $i++) {
    $newRow = $lastRow + $i

    # Column A: next day's date serial, continuing the existing sequence
    $ws.Cells.Item($newRow, 1).Value = $lastDate + $i

    # Columns B:lastCol: same values as the previous last row
    for ($c = 2; $c -le $lastCol; $c++) {
        $ws.Cells.Item($newRow, $c).Value = $rowValues[$c - 2]
    }

    # Copy the formatting (style) of the prior row's date cell onto the new one
    $ws.Cells.Item($lastRow, 1).Copy()
    $ws.Cells.Item($newRow, 1).PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = $false
